# Applies the 3.5.1.xlsx "Exakt" sheet metadata-block rework:
#  - insert 2 new rows (Geographical Area / Unit of measurement) between the
#    existing "Source" and "Comment" rows
#  - relabel/retext the surrounding rows to match the new metadata block
#  - grow the merged footnote row and bump its height
#  - drop the "Year" category-axis title from the Exakt chart
#  - remove the small second logo picture from the Exakt drawing
#  - update sheet selection/dimension bookkeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exakt")

# --- Insert two fresh rows right after the current "Source" row (36) so the
#     old Source/Copyright/footnote rows shift from 36-38 down to 38-40.
$ws.Rows("36:37").Insert()

# --- Row 35: "Specification:" / "Estimated data"  ->  "Source:" / DBDD name
$ws.Range("A35").Value = "Source:"
$ws.Range("B35").Value = "German Monitoring Center for Drugs and Drug Addiction"

# --- Row 36 (new): "Geographical Area:" / "Germany"
$ws.Range("B36").Value = "Germany"
$ws.Range("A36").Value = "Geographical Area:"

# --- Row 37 (new): "Unit of measurement:" / "Number"
$ws.Range("A37").Value = "Unit of measurement:"
$ws.Range("B37").Value = "Number"

# --- Row 38 (was the old "Source:" / DBDD row): "Comment:" / "Estimated data."
$ws.Range("A38").Value = "Comment:"
$ws.Range("B38").Value = "Estimated data."

# --- Row 39 (was the old Copyright row): keep "Copyright: " label, refresh text
$ws.Range("A39").Value = "Copyright: "
$ws.Range("B39").Value = [char]0x00A9 + " Federal Statistical Office (Destatis) 2021"

# --- Row 40 (the merged footnote row): grow its height a bit
$ws.Rows(40).RowHeight = 36.75

# --- Selection / scroll bookkeeping to match the saved view
$ws.Range("J40").Select()

# --- Widen column A so the new longer labels fit (Excel "best fit" width)
$ws.Columns("A").ColumnWidth = 17.85546875

# --- Drop the category-axis ("Year") title from the Exakt area chart
$chart = $ws.ChartObjects(1).Chart
$chart.Axes(1).HasTitle = $false

# --- Remove the small second logo picture anchored near the old row 36
$ws.Shapes.Item("Grafik 3").Delete()
